$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2473309608540925
$ws1.Range("C2").Value = 0.06208425720620843
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1169102296450939
$ws1.Range("F2").Value = 0.2486678507992895
$ws1.Range("G2").Value = 0.63249348392702
$ws1.Range("H2").Value = 0.5951043338683788
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 423
$ws1.Range("K2").Value = 111
$ws1.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2078651685393259
$ws2.Range("D2").Value = 0.3441860465116279

$ws2.Range("B3").Value = 0.06208425720620843
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1169102296450939

$ws2.Range("B4").Value = 0.2473309608540925
$ws2.Range("C4").Value = 0.2473309608540925
$ws2.Range("D4").Value = 0.2473309608540925
$ws2.Range("E4").Value = 0.2473309608540925

$ws2.Range("B5").Value = 0.5310421286031042
$ws2.Range("C5").Value = 0.6039325842696629
$ws2.Range("D5").Value = 0.2305481380783609

$ws2.Range("B6").Value = 0.9532711017825157
$ws2.Range("C6").Value = 0.2473309608540925
$ws2.Range("D6").Value = 0.3328626962051102

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 111
$ws3.Range("C2").Value = 423
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
